$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Ref, $Val)
    $rng = $Sheet.Range($Ref)
    $prevStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $prevStyle
}

$ws.Range("D2").Value = '27.563.53'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '1.662.02'
$ws.Range("E3").Value = '  -3.62%  '
$ws.Range("E4").Value = '  +0.61%  '
$ws.Range("D5").Value = '214.56'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = '0.512'
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = '23.38'
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("D9").Value = '0.259'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").Value = '1.896.28'
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("D13").Value = '1.659.56'
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("E14").Value = '  -3.13%  '
$ws.Range("D15").Value = '0.547'
$ws.Range("E15").Value = '  -3.41%  '
Set-TextValue $ws 'D16' '65.80'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").Value = '246.79'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").Value = '27.562.31'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '0.0₃0730'
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("D20").Value = '7.47'
$ws.Range("E20").Value = '  -7.29%  '
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("E23").Value = '  -4.47%  '
$ws.Range("E24").Value = '  -4.69%  '
$ws.Range("D25").Value = '145.93'
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = '7.17'
$ws.Range("E26").Value = '  -5.12%  '
Set-TextValue $ws 'D27' '16.20'
$ws.Range("E27").Value = '  -2.73%  '
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("D29").Value = '0.112'
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("E30").Value = '  +4.66%  '
$ws.Range("E31").Value = '  -1.45%  '
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").Value = '1.442.96'
$ws.Range("E33").Value = '  -6.65%  '
$ws.Range("D34").Value = '3.12'
$ws.Range("E34").Value = '  -5.88%  '
$ws.Range("E35").Value = '  -8.57%  '
$ws.Range("E36").Value = '  -0.56%  '
Set-TextValue $ws 'D37' '0.930'
$ws.Range("E37").Value = '  -4.21%  '
$ws.Range("E38").Value = '  -6.02%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("E40").Value = '  -2.80%  '
$ws.Range("D41").Value = '68.99'
$ws.Range("E41").Value = '  -3.58%  '
$ws.Range("E42").Value = '  +0.62%  '
Set-TextValue $ws 'D43' '5.40'
$ws.Range("E43").Value = '  -7.94%  '
$ws.Range("D44").Value = '0.791'
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.805.39'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("B46").Value = 'MXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D46' '2.20'
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '88.53'
$ws.Range("E48").Value = '  -4.10%  '
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("E50").Value = '  -4.34%  '
$ws.Range("D51").Value = '7.83'
$ws.Range("E51").Value = '  -6.39%  '
